$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (typo fixes / removing "(no resuelta)" suffixes) ---
$ws.Range("G4").Value = "Recuerda que si el cliente despues de tres intentos no contesta agregar esa información"
$ws.Range("G5").Value = "En plan de riesgos falta efecto de la causa no conexión por admiadmin impide conectar a maquina cliente por ejemplo."
$ws.Range("G6").Value = "Al no requerir implementacion no debe decir los mensajes de implementacion"

# --- Row 7: status becomes "Cerrada" and comment updated (no resuelta) removed ---
$ws.Range("F7").Value = "Cerrada"
$ws.Range("G7").Value = "En la parte equipo de empresa poner SOS Software y en cliente el nombre de cliente"

# --- Row 8: collapsed into a short "cd " comment row, other fields cleared ---
$ws.Range("A8").Value = "cd "
$ws.Range("B8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""

# --- Row heights ---
$ws.Range("A4").RowHeight = 55.2
$ws.Range("A5").RowHeight = 55.2
$ws.Range("A7").RowHeight = 41.75
$ws.Range("A8").RowHeight = 13.8

# --- Selection ---
$ws.Range("D3").Select()
